$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Variables")

# Update the VariableValue column (B) with the tested deployment values.
# Order matters for shared-string table layout, so write in this sequence.
$ws.Range("B8").Value = "vpc-14ee357d"
$ws.Range("B10").Value = "212.139.37.214"
$ws.Range("B11").Value = "Ryan"
$ws.Range("B13").Value = "Ryan Froggatt"
$ws.Range("B4").Value = "eu-west-2"
$ws.Range("B9").Value = "subnet-6831fd13"

$ws.Range("B5").Value = $true
$ws.Range("B6").Value = $true

# Update selected cell to reflect where the user left the cursor
$ws.Activate()
$ws.Range("B9").Select()
